# Fill in the server data row for the LoginServer Property sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (existing): A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
# Populate row 2 with the first server's data (order matches shared-string
# table append order: IP, then the server name, then the numeric-looking ID).
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "LoginServer_1"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "LoginServer_1"
$ws.Range("B2").Value = "000106001"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 6001

# Move the active selection to G5 (single cell), matching the saved view state.
$ws.Range("G5").Select() | Out-Null

# The "allow blank only" TRUE/FALSE list validation on column F now starts at
# F3 instead of F2, since F2 holds the literal IP address for row 2.
$ws.Range("F2:F1048576").Validation.Delete() | Out-Null
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"') | Out-Null
